# INTERACTIVE MAP 0.3.0 change filter for date to range
#
# The "asociado" column (F) linked each modeled/observed series to its
# counterpart series so the map could filter by an associated variable.
# Update the association codes for the "Cordillera" rows and match the
# text/number formatting already used by the "codigo_variable" column
# (B3/B4), which is what the new values represent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("Caudal Observado Cordillera"): associate with the modeled series.
$ws.Range("F3").Value = "caudales_modelados"

# Row 4 ("Caudal Modelado Cordillera"): no longer associated with a
# specific observed-cordillera code (that variable no longer exists), so
# reset it to the same placeholder used elsewhere ("-").
$ws.Range("F4").Value = "-"

# Match the formatting of the codigo_variable column for these two cells.
$ws.Range("B3").Copy()
$ws.Range("F3").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("F4").PasteSpecial(-4122)
